$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap "Periodo Mora" values between the two worker rows (1808 <-> 1807)
$ws.Range("E16").Value = "1807"
$ws.Range("E17").Value = "1808"

# Update "Salario Basico" amounts for both rows
$ws.Range("G16").Value = 1423500
$ws.Range("G17").Value = 1423500
